$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '51.747.27'
$ws.Range('E2').Value = '  +1.38%  '

$ws.Range('D3').Value = '3.031.26'
$ws.Range('E3').Value = '  +2.40%  '

$ws.Range('E4').Value = '  +0.09%  '

Set-TextValue $ws.Range('D5') '381.02'
$ws.Range('E5').Value = '  +0.49%  '

Set-TextValue $ws.Range('D6') '103.08'
$ws.Range('E6').Value = '  +1.02%  '

$ws.Range('E7').Value = '  +0.71%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('E9').Value = '  +1.53%  '

Set-TextValue $ws.Range('D10') '36.90'
$ws.Range('E10').Value = '  +1.45%  '

$ws.Range('E12').Value = '  +1.13%  '

$ws.Range('D13').Value = '3.511.17'
$ws.Range('E13').Value = '  +2.56%  '

Set-TextValue $ws.Range('D14') '18.56'
$ws.Range('E14').Value = '  +0.68%  '

$ws.Range('E15').Value = '  -0.59%  '

$ws.Range('D16').Value = '3.030.26'
$ws.Range('E16').Value = '  +2.34%  '

Set-TextValue $ws.Range('D17') '0.979'
$ws.Range('E17').Value = '  -4.14%  '

Set-TextValue $ws.Range('D18') '10.54'

$ws.Range('D19').Value = '51.741.28'
$ws.Range('E19').Value = '  +1.42%  '

Set-TextValue $ws.Range('D20') '3.08'
$ws.Range('E20').Value = '  -0.42%  '

Set-TextValue $ws.Range('D21') '12.52'
$ws.Range('E21').Value = '  +0.73%  '

$ws.Range('E22').Value = '  +1.03%  '

Set-TextValue $ws.Range('D23') '70.17'
$ws.Range('E23').Value = '  +0.78%  '

Set-TextValue $ws.Range('D24') '268.80'
$ws.Range('E24').Value = '  +0.80%  '

$ws.Range('E25').Value = '  -4.37%  '

Set-TextValue $ws.Range('D26') '8.25'
$ws.Range('E26').Value = '  +1.62%  '

Set-TextValue $ws.Range('D27') '7.65'
$ws.Range('E27').Value = '  +9.13%  '

$ws.Range('E28').Value = '  +5.24%  '

$ws.Range('E29').Value = '  -0.07%  '

Set-TextValue $ws.Range('D30') '26.26'
$ws.Range('E30').Value = '  +1.95%  '

Set-TextValue $ws.Range('D31') '0.108'
$ws.Range('E31').Value = '  +0.43%  '

Set-TextValue $ws.Range('D32') '10.30'
$ws.Range('E32').Value = '  +0.32%  '

$ws.Range('E33').Value = '  +2.42%  '

Set-TextValue $ws.Range('D34') '34.12'
$ws.Range('E34').Value = '  +0.45%  '

Set-TextValue $ws.Range('D35') '50.51'
$ws.Range('E35').Value = '  -0.07%  '

$ws.Range('E36').Value = '  +3.42%  '

$ws.Range('E38').Value = '  +4.95%  '

Set-TextValue $ws.Range('D39') '0.297'
$ws.Range('E39').Value = '  +14.40%  '

Set-TextValue $ws.Range('D40') '17.07'
$ws.Range('E40').Value = '  +2.75%  '

$ws.Range('E41').Value = '  +2.31%  '

$ws.Range('E42').Value = '  +2.48%  '

Set-TextValue $ws.Range('D43') '127.64'
$ws.Range('E43').Value = '  +6.04%  '

$ws.Range('E44').Value = '  -0.43%  '

Set-TextValue $ws.Range('D45') '3.79'
$ws.Range('E45').Value = '  +6.49%  '

Set-TextValue $ws.Range('D46') '21.81'
$ws.Range('E46').Value = '  +1.64%  '

Set-TextValue $ws.Range('D47') '2.10'
$ws.Range('E47').Value = '  +3.95%  '

$ws.Range('E48').Value = '  +4.89%  '

$ws.Range('D49').Value = '2.033.98'
$ws.Range('E49').Value = '  +1.13%  '

$ws.Range('D50').Value = '3.332.63'
$ws.Range('E50').Value = '  +2.60%  '

$ws.Range('E51').Value = '  -0.31%  '
